# Applies the "Add files via upload" commit: renames the dataset's header
# row (B1:N1, skipping the unchanged A1/K1/O1 headers) to the new short
# machine-friendly column names, and shrinks row 1's height now that the
# headers are single short words instead of wrapped multi-line text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "srcmsgformat"
$ws.Range("C1").Value = "srcprotocol"
$ws.Range("D1").Value = "int-dataformat"
$ws.Range("E1").Value = "targetmsgformat"
$ws.Range("F1").Value = "targetprotocol"
$ws.Range("G1").Value = "Interface-type"
$ws.Range("H1").Value = "msgfieldcount"
$ws.Range("I1").Value = "rulecount"
$ws.Range("J1").Value = "operationcount"
$ws.Range("L1").Value = "product"
$ws.Range("M1").Value = "new-existing"
$ws.Range("N1").Value = "exposed-as-api"

# Header text is now a single short word per cell instead of a wrapped
# multi-line label, so the row no longer needs to be as tall.
$ws.Rows.Item(1).RowHeight = 45
